# changing names and moving append to writer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 11 "manager" records that get duplicated (originally rows 1-11)
$records = @(
    @(3239632203, "Tu", "Dagny", 1386825528, 8046993693),
    @(8795146687, "Kecia", "Callie", 1386825528, 8046993693),
    @(7569690984, "Kalyn", "Leonardo", 1386825528, 8046993693),
    @(317440315, "Evelynn", "Anamaria", 1386825528, 8046993693),
    @(3502441984, "Harvey", "Ty", 7285813456, 8046993693),
    @(6841158304, "Nakesha", "Chana", 7285813456, 8046993693),
    @(3273508974, "Maegan", "Josie", 7285813456, 8046993693),
    @(3788855002, "Wendie", "Elias", 7285813456, 8046993693),
    @(5002807220, "Keven", "Vivian", 7285813456, 8046993693),
    @(3569528995, "Marquetta", "Darla", 7285813456, 8046993693),
    @(3723735449, "Dorothea", "Deloris", 7285813456, 8046993693)
)

# The final row (previously row 12) that now moves to row 34
$lastRow = @(9372070413, "Dillon", "Jayne", 7285813456, 8046993693)

# Append the 11 records twice starting at row 12 (rows 12-22, then 23-33)
$r = 12
for ($pass = 0; $pass -lt 2; $pass++) {
    foreach ($row in $records) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $r = $r + 1
    }
}

# Write the moved last row at row 34
$ws.Cells.Item($r, 1).Value = $lastRow[0]
$ws.Cells.Item($r, 2).Value = $lastRow[1]
$ws.Cells.Item($r, 3).Value = $lastRow[2]
$ws.Cells.Item($r, 4).Value = $lastRow[3]
$ws.Cells.Item($r, 5).Value = $lastRow[4]
